# Apply cell-value updates produced by the scheduled Sheets runner.
# All target cells are plain numeric (t="n") cells with cached values -
# no formulas are involved, so this is a straightforward Value2 rewrite
# per (sheet, row, column), including clearing cells that become empty
# and populating cells that previously had no value at all.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 68
$ws.Range("H68").Value2 = 30000
$ws.Range("J68").Value2 = 30000
$ws.Range("L68").Value2 = 30000
$ws.Range("N68").Value2 = -31498
# Row 71
$ws.Range("H71").Value2 = 30000
$ws.Range("J71").Value2 = 30000
$ws.Range("L71").Value2 = 90000
$ws.Range("N71").Value2 = -97488
# Row 138
$ws.Range("H138").Value2 = 4275877
$ws.Range("I138").Value2 = 7937986
$ws.Range("J138").Value2 = 3415.9167
$ws.Range("K138").Value2 = 23813958
$ws.Range("L138").Value2 = 10247.7501
$ws.Range("M138").Value2 = -23808818
$ws.Range("N138").Value2 = -20527.7501

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Range("H22").Value2 = 238.66667
$ws.Range("I22").Value2 = 238.66667
$ws.Range("K22").Value2 = 238.66667
$ws.Range("M22").Value2 = 60.33332999999999
# Row 32
$ws.Range("H32").Value2 = 2890
$ws.Range("I32").Value2 = 2860.606
$ws.Range("J32").Value2 = 5800
$ws.Range("K32").Value2 = 2860.606
$ws.Range("L32").Value2 = 5800
$ws.Range("M32").Value2 = -2573.606
$ws.Range("N32").Value2 = -6374
# Row 34
$ws.Range("H34").Value2 = 16018.667
$ws.Range("J34").Value2 = 16018.667
$ws.Range("L34").Value2 = 16018.667
$ws.Range("N34").Value2 = -16560.667
# Row 38
$ws.Range("H38").Value2 = 11000
$ws.Range("I38").Value2 = 2000
$ws.Range("K38").Value2 = 2000
$ws.Range("M38").Value2 = -1533
# Row 61
$ws.Range("H61").Value2 = 3239.682
$ws.Range("I61").Value2 = 3345.6667
$ws.Range("J61").Value2 = 1014
$ws.Range("K61").Value2 = 3345.6667
$ws.Range("L61").Value2 = 1014
$ws.Range("M61").Value2 = -3133.6667
$ws.Range("N61").Value2 = -1438
# Row 102
$ws.Range("H102").Value2 = 770
$ws.Range("I102").Value2 = 770
$ws.Range("J102").Value2 = 0
$ws.Range("K102").Value2 = 770
$ws.Range("L102").Value2 = 0
$ws.Range("M102").Value2 = 852
$ws.Range("N102").ClearContents()
# Row 136
$ws.Range("H136").Value2 = 3239.682
$ws.Range("I136").Value2 = 3345.6667
$ws.Range("J136").Value2 = 1014
$ws.Range("K136").Value2 = 10037.0001
$ws.Range("L136").Value2 = 3042
$ws.Range("M136").Value2 = -7487.000100000001
$ws.Range("N136").Value2 = -8142

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 49
$ws.Range("H49").Value2 = 0
$ws.Range("J49").Value2 = 0
$ws.Range("L49").Value2 = 0
$ws.Range("N49").ClearContents()
# Row 63
$ws.Range("H63").Value2 = 42771
$ws.Range("J63").Value2 = 42771
$ws.Range("L63").Value2 = 42771
$ws.Range("N63").Value2 = -44143
# Row 66
$ws.Range("H66").Value2 = 42771
$ws.Range("J66").Value2 = 42771
$ws.Range("L66").Value2 = 128313
$ws.Range("N66").Value2 = -135177
# Row 134
$ws.Range("H134").Value2 = 3840.8845
$ws.Range("I134").Value2 = 2677.28
$ws.Range("J134").Value2 = 4918.2964
$ws.Range("K134").Value2 = 8031.84
$ws.Range("L134").Value2 = 14754.8892
$ws.Range("M134").Value2 = -5496.84
$ws.Range("N134").Value2 = -19824.8892

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value2 = 1831.08
$ws.Range("I99").Value2 = 1612.1428
$ws.Range("J99").Value2 = 2109.7273
$ws.Range("K99").Value2 = 1612.1428
$ws.Range("L99").Value2 = 2109.7273
$ws.Range("M99").Value2 = -114.1428000000001
$ws.Range("N99").Value2 = -5105.7273
# Row 126
$ws.Range("H126").Value2 = 1831.08
$ws.Range("I126").Value2 = 1612.1428
$ws.Range("J126").Value2 = 2109.7273
$ws.Range("K126").Value2 = 4836.428400000001
$ws.Range("L126").Value2 = 6329.1819
$ws.Range("M126").Value2 = -2366.428400000001
$ws.Range("N126").Value2 = -11269.1819

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value2 = 634.9355
$ws.Range("I5").Value2 = 219.29411
$ws.Range("J5").Value2 = 1139.6428
$ws.Range("K5").Value2 = 657.8823299999999
$ws.Range("L5").Value2 = 3418.9284
$ws.Range("M5").Value2 = -545.8823299999999
$ws.Range("N5").Value2 = -3642.9284
# Row 109
$ws.Range("H109").Value2 = 3411.5
$ws.Range("I109").Value2 = 912.2
$ws.Range("J109").Value2 = 4800
$ws.Range("K109").Value2 = 2736.6
$ws.Range("L109").Value2 = 14400
$ws.Range("M109").Value2 = -1696.6
$ws.Range("N109").Value2 = -16480
# Row 135
$ws.Range("H135").Value2 = 634.9355
$ws.Range("I135").Value2 = 219.29411
$ws.Range("J135").Value2 = 1139.6428
$ws.Range("K135").Value2 = 1973.64699
$ws.Range("L135").Value2 = 10256.7852
$ws.Range("M135").Value2 = 561.35301
$ws.Range("N135").Value2 = -15326.7852

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value2 = 36666.668
# Row 73
$ws.Range("H73").Value2 = 36666.668
# Row 102
$ws.Range("H102").Value2 = 3485.2
$ws.Range("I102").Value2 = 3863.2632
$ws.Range("K102").Value2 = 3863.2632
$ws.Range("M102").Value2 = -2241.2632
# Row 132
$ws.Range("H132").Value2 = 4371.8184
$ws.Range("I132").Value2 = 4598.5557
$ws.Range("J132").Value2 = 3351.5
$ws.Range("K132").Value2 = 13795.6671
$ws.Range("L132").Value2 = 10054.5
$ws.Range("M132").Value2 = -11265.6671
$ws.Range("N132").Value2 = -15114.5
# Row 136
$ws.Range("H136").Value2 = 30755.143
$ws.Range("I136").Value2 = 49326
$ws.Range("J136").Value2 = 27660
$ws.Range("K136").Value2 = 147978
$ws.Range("L136").Value2 = 82980
$ws.Range("M136").Value2 = -145428
$ws.Range("N136").Value2 = -88080

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value2 = 7084.2104
$ws.Range("I7").Value2 = 10833.333
$ws.Range("K7").Value2 = 10833.333
$ws.Range("M7").Value2 = -10721.333
# Row 24
$ws.Range("H24").Value2 = 0
$ws.Range("J24").Value2 = 0
$ws.Range("L24").Value2 = 0
$ws.Range("N24").ClearContents()
# Row 86
$ws.Range("H86").Value2 = 33695
$ws.Range("J86").Value2 = 33695
$ws.Range("L86").Value2 = 33695
$ws.Range("N86").Value2 = -36067
# Row 87
$ws.Range("H87").Value2 = 26353.334
$ws.Range("I87").Value2 = 10171
$ws.Range("J87").Value2 = 34444.5
$ws.Range("K87").Value2 = 10171
$ws.Range("L87").Value2 = 34444.5
$ws.Range("M87").Value2 = -9048
$ws.Range("N87").Value2 = -36690.5
# Row 89
$ws.Range("H89").Value2 = 33695
$ws.Range("J89").Value2 = 33695
$ws.Range("L89").Value2 = 101085
$ws.Range("N89").Value2 = -112941
# Row 90
$ws.Range("H90").Value2 = 26353.334
$ws.Range("I90").Value2 = 10171
$ws.Range("J90").Value2 = 34444.5
$ws.Range("K90").Value2 = 30513
$ws.Range("L90").Value2 = 103333.5
$ws.Range("M90").Value2 = -24897
$ws.Range("N90").Value2 = -114565.5
# Row 126
$ws.Range("H126").Value2 = 7084.2104
$ws.Range("I126").Value2 = 10833.333
$ws.Range("K126").Value2 = 32499.999
$ws.Range("M126").Value2 = -30029.999
# Row 136
$ws.Range("H136").Value2 = 3638.451
$ws.Range("I136").Value2 = 1719.9166
$ws.Range("J136").Value2 = 34335
$ws.Range("K136").Value2 = 5159.7498
$ws.Range("L136").Value2 = 103005
$ws.Range("M136").Value2 = -2609.7498
$ws.Range("N136").Value2 = -108105

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value2 = 4461.5
$ws.Range("I62").Value2 = 4346.3076
$ws.Range("J62").Value2 = 4627.8887
$ws.Range("K62").Value2 = 4346.3076
$ws.Range("L62").Value2 = 4627.8887
$ws.Range("M62").Value2 = -3722.3076
$ws.Range("N62").Value2 = -5875.8887
# Row 65
$ws.Range("H65").Value2 = 4461.5
$ws.Range("I65").Value2 = 4346.3076
$ws.Range("J65").Value2 = 4627.8887
$ws.Range("K65").Value2 = 21731.538
$ws.Range("L65").Value2 = 23139.4435
$ws.Range("M65").Value2 = -18611.538
$ws.Range("N65").Value2 = -29379.4435
# Row 75
$ws.Range("H75").Value2 = 39400
$ws.Range("J75").Value2 = 39400
$ws.Range("L75").Value2 = 39400
$ws.Range("N75").Value2 = -41272
# Row 78
$ws.Range("H78").Value2 = 39400
$ws.Range("J78").Value2 = 39400
$ws.Range("L78").Value2 = 118200
$ws.Range("N78").Value2 = -127560
# Row 80
$ws.Range("H80").Value2 = 40300.5
$ws.Range("J80").Value2 = 40300.5
$ws.Range("L80").Value2 = 40300.5
$ws.Range("N80").Value2 = -42296.5
# Row 83
$ws.Range("H83").Value2 = 40300.5
$ws.Range("J83").Value2 = 40300.5
$ws.Range("L83").Value2 = 120901.5
$ws.Range("N83").Value2 = -130885.5
